$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Color constant matching the existing "blue" (FF0070C0) font style already
# used elsewhere in the sheet (e.g. C46), expressed as the BGR long Excel
# COM expects for Font.Color (0x00BBGGRR).
$blue = 12611584

# --- Row 46 ("Option to toggle between short videos/burst"): mark the
# D/E columns as OK (previously "Need test"), matching the row's styling.
$ws.Range("D46").Value = "OK"
$ws.Range("D46").Font.Color = $blue

$ws.Range("E46").Value = "OK"
$ws.Range("E46").Font.Color = $blue

# --- Row 47 (new): "Astrophotography always-on button" feature row, all
# columns marked OK with the same blue styling used throughout the table.
$ws.Range("A47").Value = "Astrophotography always-on button"
$ws.Range("A47").Font.Color = $blue

$ws.Range("B47").Value = "OK"
$ws.Range("B47").Font.Color = $blue

$ws.Range("C47").Value = "OK"
$ws.Range("C47").Font.Color = $blue

$ws.Range("D47").Value = "OK"
$ws.Range("D47").Font.Color = $blue

$ws.Range("E47").Value = "OK"
$ws.Range("E47").Font.Color = $blue

$ws.Range("F47").Value = "OK"
$ws.Range("F47").Font.Color = $blue

# Match the row height used by the other compact rows (e.g. row 46).
$ws.Rows.Item(47).RowHeight = 12.8

# Update the active selection to reflect where editing ended up.
$ws.Range("H46").Select()
